# calorimetry : scripts : data load : all setup info accumulated into setup file/sheet
#
# The "targets" worksheet is removed; the header pair it used to hold
# ("constants " / "Comp") is now appended as a new row on the "setup" sheet.
# The previously-second-to-last sheet ("input_stoich_coefficients") becomes
# the active sheet afterwards (it slides into the tab position vacated by
# "targets").

$wb = $excel.ActiveWorkbook

# 1) Append the "constants " / "Comp" row to the "setup" sheet.
$setup = $wb.Worksheets.Item("setup")
$setup.Range("A4").Value = "constants "
$setup.Range("B4").Value = "Comp"

# 2) Remove the now-redundant "targets" sheet.
$targets = $wb.Worksheets.Item("targets")
$null = $targets.Delete()

# 3) Move the selection on "setup" up to A3 (matches the recorded view state).
$null = $setup.Range("A3").Select()

# 4) Activate "input_stoich_coefficients" (now occupying the slot after
#    "heats") and restore its own recorded selection.
$stoich = $wb.Worksheets.Item("input_stoich_coefficients")
$null = $stoich.Activate()
$null = $stoich.Range("K6").Select()
